$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 (CONFIDENCE.T label, uses CONFIDENCE.NORM formula - matches source data exactly,
# including the original commit's copy/paste quirk) must be written BEFORE row 9 so the
# shared-string table gets "CONFIDENCE.T" at index 9 and "CONFIDENCE.NORM" at index 10,
# matching the target workbook.
$ws.Range("A10").Value = "CONFIDENCE.T"
$ws.Range("B10").Formula = "=_xlfn.CONFIDENCE.NORM(C10,D10,E10)"
$ws.Range("C10").Formula = "=2/15"
$ws.Range("D10").Value = 6.6
$ws.Range("E10").Value = 44

$ws.Range("A9").Value = "CONFIDENCE.NORM"
$ws.Range("B9").Formula = "=_xlfn.CONFIDENCE.NORM(C9,D9,E9)"
$ws.Range("C9").Formula = "=2/15"
$ws.Range("D9").Value = 6.6
$ws.Range("E9").Value = 44

# Column A widened to fit the new, longer function names.
$ws.Columns.Item(1).ColumnWidth = 18.5

# New selection left behind by the editor, one row below the new data.
[void]$ws.Range("A11").Select()
